$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - match styling of existing header cells (bold, centered, bordered)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Data cells F2:F12 - plain text timestamps
$ws.Range("F2").Value = "2021-10-05 13:41:02.859198"
$ws.Range("F3").Value = "2021-10-05 13:41:02.859208"
$ws.Range("F4").Value = "2021-10-05 13:41:02.859210"
$ws.Range("F5").Value = "2021-10-05 13:41:02.859212"
$ws.Range("F6").Value = "2021-10-05 13:41:02.859215"
$ws.Range("F7").Value = "2021-10-05 13:41:02.859217"
$ws.Range("F8").Value = "2021-10-05 13:41:02.859219"
$ws.Range("F9").Value = "2021-10-05 13:41:02.859221"
$ws.Range("F10").Value = "2021-10-05 13:41:02.859223"
$ws.Range("F11").Value = "2021-10-05 13:41:02.859225"
$ws.Range("F12").Value = "2021-10-05 13:41:02.859227"
